# --- Part 1: insert the new "Meta description" paragraph right after the title ---
$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our expert review of Drago: Jewels of Fortune and play for free. Experience the dragon-themed design, streak respin feature, and free spins with multipliers.</w:t></w:r><w:r/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($metaXml) | Out-Null

# InsertXML split the title paragraph's tail into an extra blank paragraph
# (so the new content becomes its own paragraph) - remove that blank spacer.
$blank = $d.Paragraphs(3)
$blank.Range.Delete() | Out-Null

# --- Part 2: drop the trailing duplicate title paragraph, and turn the ---
# --- trailing italic description into the image-prompt text instead.  ---
$n = $d.Paragraphs.Count
$titleDup = $d.Paragraphs($n - 1)
$titleDup.Range.Delete() | Out-Null

$n2 = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($n2)
$fullRange = $italicPara.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.Text = "Please create an image for ""Drago: Jewels of Fortune"" featuring a happy Maya warrior with glasses in a cartoon style. The warrior should be depicted in a jungle setting with treasure chests and dragons in the background. The image should be vibrant and colorful, with attention to detail in the warrior's clothing and accessories. The overall vibe should be adventurous and exciting, reflecting the theme of the game. Thank you!"
